$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 78, pushing existing rows 78-128 down to 80-130.
$ws.Rows.Item(78).Insert()
$ws.Rows.Item(78).Insert()

# New row 78: Chirimoya "Especial", 2022-11-14 (serial 44879)
$ws.Cells.Item(78,1).Value = 5
$ws.Cells.Item(78,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(78,3).Value = "Maule"
$ws.Cells.Item(78,4).Value = 44879
$ws.Cells.Item(78,5).Value = 7
$ws.Cells.Item(78,6).Value = "Fruta"
$ws.Cells.Item(78,7).Value = 100107
$ws.Cells.Item(78,8).Value = "Otros"
$ws.Cells.Item(78,9).Value = 100107002
$ws.Cells.Item(78,10).Value = "Chirimoya"
$ws.Cells.Item(78,11).Value = "Cultivar IV Región"
$ws.Cells.Item(78,12).Value = "Especial"
$ws.Cells.Item(78,13).Value = 190
$ws.Cells.Item(78,14).Value = 25000
$ws.Cells.Item(78,15).Value = 25000
$ws.Cells.Item(78,16).Value = 25000
$ws.Cells.Item(78,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(78,18).Value = "Provincia de Limarí"
$ws.Cells.Item(78,19).Value = 2500
$ws.Cells.Item(78,20).Value = 10

# New row 79: Chirimoya "Primera", 2022-11-14 (serial 44879)
$ws.Cells.Item(79,1).Value = 5
$ws.Cells.Item(79,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(79,3).Value = "Maule"
$ws.Cells.Item(79,4).Value = 44879
$ws.Cells.Item(79,5).Value = 7
$ws.Cells.Item(79,6).Value = "Fruta"
$ws.Cells.Item(79,7).Value = 100107
$ws.Cells.Item(79,8).Value = "Otros"
$ws.Cells.Item(79,9).Value = 100107002
$ws.Cells.Item(79,10).Value = "Chirimoya"
$ws.Cells.Item(79,11).Value = "Cultivar IV Región"
$ws.Cells.Item(79,12).Value = "Primera"
$ws.Cells.Item(79,13).Value = 200
$ws.Cells.Item(79,14).Value = 23000
$ws.Cells.Item(79,15).Value = 23000
$ws.Cells.Item(79,16).Value = 23000
$ws.Cells.Item(79,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(79,18).Value = "Provincia de Limarí"
$ws.Cells.Item(79,19).Value = 2300
$ws.Cells.Item(79,20).Value = 10
